# Add a new PROPERTY_CODES sheet describing the property_codes table schema,
# matching the style/layout of the other schema-description sheets
# (LOCATIONS, sqlite_sequence, GOOGLE_SOLAR, CEJST).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (CEJST) so it lands as
# sheet #5 / last tab, with sheetId=5 / rId5.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "PROPERTY_CODES"

# Copy the header row (values + bold/centered/thin-bordered formatting) from
# an existing sheet so the new sheet reuses the same shared strings & style
# index instead of minting new ones.
$srcHeader = $wb.Worksheets.Item("CEJST").Range("A1:F1")
$srcHeader.Copy($ws.Range("A1:F1"))

# Data rows describing each column of the PROPERTY_CODES table:
#   cid, column_name, data_type, notnull, default_value, pk
$rows = @(
    @(0, "property_code_id", "INTEGER",  0, "",                  1),
    @(1, "property_code",    "TEXT",     0, "",                  0),
    @(2, "description",      "TEXT",     0, "",                  0),
    @(3, "name",             "TEXT",     0, "",                  0),
    @(4, "date_added",       "DATETIME", 0, "CURRENT_TIMESTAMP", 0)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    if ($row[4] -ne "") {
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
